# Apply the "Early version of any strategy" edit to the TestCases sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: switch exchange to ByBit, move the "From" date earlier, and use MACD as the strategy
$ws.Range("B2").Value = "ByBit"
$ws.Range("D2").Value = (Get-Date -Year 2021 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("J2").Value = "MACD"

# Row 3: move the "From" date earlier and switch the strategy to the new EarlyMACD variant
$ws.Range("D3").Value = (Get-Date -Year 2021 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("J3").Value = "EarlyMACD"

# Row 4 (the old standalone EarlyMACD test case) is no longer needed now that Early
# strategy variants can be created directly via multiple inheritance - remove its
# test entry entirely, leaving only the date-formatted D4/E4 cells behind.
$ws.Range("A4").Clear()
$ws.Range("B4").Clear()
$ws.Range("C4").Clear()
$ws.Range("F4").Clear()
$ws.Range("G4").Clear()
$ws.Range("H4").Clear()
$ws.Range("I4").Clear()
$ws.Range("J4").Clear()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Restore the cursor/selection position recorded in the saved workbook
[void]$ws.Range("M7").Select()
